$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, $row, $values) {
    $c = 1
    foreach ($v in $values) {
        $ws.Cells.Item($row, $c).Value = $v
        $c++
    }
}

function Extend-CF($ws, $col, $firstRow, $newLastRow) {
    $addr = $col + $firstRow
    $r = $ws.Range($addr)
    $cnt = $r.FormatConditions.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $fc = $r.FormatConditions.Item($i)
        $newRange = $ws.Range($col + $firstRow + ":" + $col + $newLastRow)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# ------------------------------------------------------------------
# Sheet 1 : "Test Parallel-Series System"  (A1:P6 -> A1:P10)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# -- fix existing rows 3-6 --
$ws1.Cells.Item(3,4).Value = 3    # D3: 2 -> 3

$ws1.Cells.Item(4,4).Value = 3    # D4: 0 -> 3

$ws1.Cells.Item(5,2).Value = 3    # B5: 1 -> 3
$ws1.Cells.Item(5,3).Value = 2    # C5: 1 -> 2
$ws1.Cells.Item(5,4).Value = 3    # D5: 0 -> 3

$ws1.Cells.Item(6,2).Value = 3    # B6: 0 -> 3
$ws1.Cells.Item(6,3).Value = 2    # C6: 0 -> 2
$ws1.Cells.Item(6,4).Value = 3    # D6: 0 -> 3

# -- add new rows 7-10, cloning formatting from row 6 --
$ws1.Range("A6:P6").Copy()
$ws1.Range("A7:P10").PasteSpecial(-4122)

Set-RowValues $ws1 7  @(5,3,2,3,3,3,3,3,3,3,3,3,3,3,3)
Set-RowValues $ws1 8  @(6,3,2,3,3,3,3,3,3,3,3,3,3,3,3)
Set-RowValues $ws1 9  @(7,3,0,3,3,3,2,3,3,3,3,3,3,3,3)
Set-RowValues $ws1 10 @(8,0,0,3,3,2,2,0,3,3,3,3,3,3,3)

$ws1.Cells.Item(7,16).Formula  = "=IF(B7 = I7, 1, 0)"
$ws1.Cells.Item(8,16).Formula  = "=IF(B8 = I8, 1, 0)"
$ws1.Cells.Item(9,16).Formula  = "=IF(B9 = I9, 1, 0)"
$ws1.Cells.Item(10,16).Formula = "=IF(B10 = I10, 1, 0)"

# -- extend conditional formatting ranges --
Extend-CF $ws1 "C" 2 11
Extend-CF $ws1 "D" 2 11
Extend-CF $ws1 "F" 2 11
Extend-CF $ws1 "G" 2 11
Extend-CF $ws1 "J" 2 11
Extend-CF $ws1 "K" 2 11
Extend-CF $ws1 "M" 2 11
Extend-CF $ws1 "N" 2 11
Extend-CF $ws1 "P" 2 10

# ------------------------------------------------------------------
# Sheet 2 : "Sensed comp_1 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(3,3).Value = 1    # C3: 0 -> 1

$ws2.Cells.Item(5,2).Value = 2    # B5: 1 -> 2
$ws2.Cells.Item(5,9).Value = 2    # I5: 1 -> 2

$ws2.Cells.Item(6,2).Value = 2    # B6: 0 -> 2
$ws2.Cells.Item(6,9).Value = 2    # I6: 0 -> 2

$ws2.Range("A6:K6").Copy()
$ws2.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws2 7  @(5,2,0,0,1,3,3,3,2)
Set-RowValues $ws2 8  @(6,2,0,0,1,3,3,3,2)
Set-RowValues $ws2 9  @(7,0,0,0,0,3,3,3,2)
Set-RowValues $ws2 10 @(8,0,0,0,0,3,3,3,2)

for ($r = 7; $r -le 10; $r++) {
    $ws2.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws2.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws2 "J" 2 10
Extend-CF $ws2 "K" 2 10

# ------------------------------------------------------------------
# Sheet 3 : "Sensed comp_2 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(3,2).Value = 3    # B3: 2 -> 3
$ws3.Cells.Item(3,5).Value = 1    # E3: 0 -> 1

$ws3.Cells.Item(4,2).Value = 3    # B4: 0 -> 3
$ws3.Cells.Item(5,2).Value = 3    # B5: 0 -> 3
$ws3.Cells.Item(6,2).Value = 3    # B6: 0 -> 3

$ws3.Range("A6:K6").Copy()
$ws3.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws3 7  @(5,3,0,0,0,3,3,3,3)
Set-RowValues $ws3 8  @(6,3,0,0,0,3,3,3,3)
Set-RowValues $ws3 9  @(7,3,0,0,0,3,3,3,3)
Set-RowValues $ws3 10 @(8,3,0,0,0,3,3,3,3)

for ($r = 7; $r -le 10; $r++) {
    $ws3.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws3.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws3 "J" 2 10
Extend-CF $ws3 "K" 2 10

# ------------------------------------------------------------------
# Sheet 4 : "Sensed comp_3 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3,4).Value = 0    # D3: 1 -> 0
$ws4.Cells.Item(4,4).Value = 0    # D4: 1 -> 0
$ws4.Cells.Item(5,4).Value = 0    # D5: 1 -> 0

$ws4.Cells.Item(6,4).Value = 0    # D6: 1 -> 0
$ws4.Cells.Item(6,5).Value = 1    # E6: 0 -> 1

$ws4.Range("A6:K6").Copy()
$ws4.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws4 7  @(5,3,0,0,1,3,3,3,3)
Set-RowValues $ws4 8  @(6,3,0,0,0,3,3,3,3)
Set-RowValues $ws4 9  @(7,3,0,0,0,3,3,3,3)
Set-RowValues $ws4 10 @(8,3,0,0,0,3,3,3,3)

for ($r = 7; $r -le 10; $r++) {
    $ws4.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws4.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws4 "J" 2 10
Extend-CF $ws4 "K" 2 10

# ------------------------------------------------------------------
# Sheet 5 : "Sensed comp_4 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Cells.Item(3,5).Value = 1    # E3: 0 -> 1
$ws5.Cells.Item(4,4).Value = 1    # D4: 0 -> 1
$ws5.Cells.Item(5,3).Value = 0    # C5: 1 -> 0
$ws5.Cells.Item(6,3).Value = 0    # C6: 1 -> 0

$ws5.Range("A6:K6").Copy()
$ws5.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws5 7  @(5,3,0,0,0,3,3,3,3)
Set-RowValues $ws5 8  @(6,3,0,0,0,3,3,3,3)
Set-RowValues $ws5 9  @(7,3,0,0,0,3,3,3,3)
Set-RowValues $ws5 10 @(8,2,0,0,0,3,3,3,3)

for ($r = 7; $r -le 10; $r++) {
    $ws5.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws5.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws5 "J" 2 10
Extend-CF $ws5 "K" 2 10

# ------------------------------------------------------------------
# Sheet 6 : "Sensed comp_5 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Cells.Item(3,3).Value = 0    # C3: 1 -> 0

$ws6.Cells.Item(4,3).Value = 0    # C4: 1 -> 0
$ws6.Cells.Item(4,5).Value = 1    # E4: 0 -> 1

$ws6.Cells.Item(5,3).Value = 0    # C5: 1 -> 0
$ws6.Cells.Item(5,4).Value = 0    # D5: 1 -> 0
$ws6.Cells.Item(5,5).Value = 1    # E5: 0 -> 1

$ws6.Cells.Item(6,3).Value = 0    # C6: 1 -> 0
$ws6.Cells.Item(6,4).Value = 0    # D6: 1 -> 0
$ws6.Cells.Item(6,5).Value = 1    # E6: 0 -> 1

$ws6.Range("A6:K6").Copy()
$ws6.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws6 7  @(5,3,0,0,0,3,3,3,3)
Set-RowValues $ws6 8  @(6,3,0,0,0,3,3,3,3)
Set-RowValues $ws6 9  @(7,2,0,0,0,3,3,3,3)
Set-RowValues $ws6 10 @(8,2,0,0,0,3,3,3,3)

for ($r = 7; $r -le 10; $r++) {
    $ws6.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws6.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws6 "J" 2 10
Extend-CF $ws6 "K" 2 10

# ------------------------------------------------------------------
# Sheet 7 : "Sensed comp_6 History"  (A1:K6 -> A1:K10)
# ------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

$ws7.Cells.Item(4,4).Value = 0    # D4: 1 -> 0

$ws7.Cells.Item(5,3).Value = 1    # C5: 0 -> 1
$ws7.Cells.Item(5,5).Value = 0    # E5: 1 -> 0

$ws7.Cells.Item(6,3).Value = 1    # C6: 0 -> 1
$ws7.Cells.Item(6,5).Value = 0    # E6: 1 -> 0

$ws7.Range("A6:K6").Copy()
$ws7.Range("A7:K10").PasteSpecial(-4122)

Set-RowValues $ws7 7  @(5,3,0,0,0,3,3,3,3)
Set-RowValues $ws7 8  @(6,3,0,0,0,3,3,3,3)
Set-RowValues $ws7 9  @(7,3,0,0,0,3,3,3,3)
Set-RowValues $ws7 10 @(8,0,0,0,0,3,3,3,3)

for ($r = 7; $r -le 10; $r++) {
    $ws7.Cells.Item($r,10).Formula = "=IF(B$r = F$r, 1, 0)"
    $ws7.Cells.Item($r,11).Formula = "=MODE(C${r}:E${r})"
}

Extend-CF $ws7 "J" 2 10
Extend-CF $ws7 "K" 2 10
